# fix "nan" in blank text cells bug
# - Special_Rules (column C) cells that literally hold the text "nan"
#   (an artifact of a pandas/NaN -> string migration) are blanked out.
# - Weapon/Armor/Shield checkbox columns (K/L/M) are mutually exclusive
#   flags driven by the item's Type (column H). Any checkbox that does not
#   correspond to the item's actual Type is blanked out instead of TRUE.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

# Data starts on row 3 (rows 1-2 hold migration metadata / column headers).
for ($row = 3; $row -le $lastRow; $row++) {

    # --- Special_Rules: blank out the literal "nan" placeholder text ---
    $specialRules = $ws.Cells.Item($row, 3)
    if ($specialRules.Value2 -eq "nan") {
        $specialRules.Value2 = ""
    }

    # --- Weapon (K) / Armor (L) / Shield (M) checkboxes ---
    $type = $ws.Cells.Item($row, 8).Value2

    $weaponCell = $ws.Cells.Item($row, 11)
    $armorCell  = $ws.Cells.Item($row, 12)
    $shieldCell = $ws.Cells.Item($row, 13)

    $isWeapon = ($type -eq "Fighting Weapon") -or ($type -eq "Ranged Weapon")
    $isArmor  = ($type -eq "Armor")
    $isShield = ($type -eq "Shield")

    if ((-not $isWeapon) -and ($weaponCell.Value2 -eq $true)) {
        $weaponCell.Value2 = ""
    }
    if ((-not $isArmor) -and ($armorCell.Value2 -eq $true)) {
        $armorCell.Value2 = ""
    }
    if ((-not $isShield) -and ($shieldCell.Value2 -eq $true)) {
        $shieldCell.Value2 = ""
    }
}
